# Insert two new rows at the top of the weekly block (before existing row 461),
# pushing all rows 461:496 down to 463:498 (this also grows the table from 496
# to 498 rows and keeps the date-formatted style on column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("461:462").Insert()

# New row 461: Femacal de La Calera / Cebollin / Primera
$ws.Cells.Item(461, 1).Value = 3
$ws.Cells.Item(461, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(461, 3).Value = "Coquimbo"
$ws.Cells.Item(461, 4).Value = 44746
$ws.Cells.Item(461, 5).Value = 5
$ws.Cells.Item(461, 6).Value = 100112037
$ws.Cells.Item(461, 7).Value = "Cebollín"
$ws.Cells.Item(461, 8).Value = "Sin especificar"
$ws.Cells.Item(461, 9).Value = "Primera"
$ws.Cells.Item(461, 10).Value = 145
$ws.Cells.Item(461, 11).Value = 6500
$ws.Cells.Item(461, 12).Value = 7000
$ws.Cells.Item(461, 13).Value = 6759
$ws.Cells.Item(461, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(461, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(461, 16).Value = 188
$ws.Cells.Item(461, 17).Value = 36
$ws.Cells.Item(461, 18).Value = "Hortaliza"

# New row 462: Femacal de La Calera / Cebollin / Segunda
$ws.Cells.Item(462, 1).Value = 3
$ws.Cells.Item(462, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(462, 3).Value = "Coquimbo"
$ws.Cells.Item(462, 4).Value = 44746
$ws.Cells.Item(462, 5).Value = 5
$ws.Cells.Item(462, 6).Value = 100112037
$ws.Cells.Item(462, 7).Value = "Cebollín"
$ws.Cells.Item(462, 8).Value = "Sin especificar"
$ws.Cells.Item(462, 9).Value = "Segunda"
$ws.Cells.Item(462, 10).Value = 67
$ws.Cells.Item(462, 11).Value = 5500
$ws.Cells.Item(462, 12).Value = 5500
$ws.Cells.Item(462, 13).Value = 5500
$ws.Cells.Item(462, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(462, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(462, 16).Value = 153
$ws.Cells.Item(462, 17).Value = 36
$ws.Cells.Item(462, 18).Value = "Hortaliza"

# Make sure column D keeps the date-time numeric format used by the rest of
# the "Fecha" column for the two freshly written cells.
$ws.Range("D461:D462").NumberFormat = "YYYY-MM-DD HH:MM:SS"
